# Daily attendance processing - 2026-02-07 09:06:19 UTC
# Update "Recorded By" (column G) values from "Administrator, Miss Dina Nasr"
# to "Miss Dina Nasr, Administrator" for the specific rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(9,10,12,14,15,16,17,18,21,23,35,36,38,40,41,42,43,44,47,49,61,62,64,66,67,68,69,70,73,75,87,88,90,92,93,94,95,96,99,101,113,114,116,118,119,120,121,122,125,127,139,140,142,144,145,146,147,148,151,153,164,167,170,174,175,176,177,191,194,197,201,202,203,204,218,221,224,228,229,230,231,245,248,251,255,256,257,258,272,275,278,282,283,284,285,299,302,305,309,310,311,312)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Miss Dina Nasr, Administrator"
}
